# Re-apply the new default table style to the table on slide 16.
# (Mirrors what PowerPoint does when a user picks a different style from
# the Table Design gallery: Shape.Table.ApplyStyle("{guid}") rewrites the
# <a:tableStyleId> element inside the table's tblPr.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

$tableShape = $null
foreach ($sh in $s.Shapes) {
    if ($sh.HasTable) {
        $tableShape = $sh
        break
    }
}

$tableShape.Table.ApplyStyle("{E81436E4-EEC1-4C42-B0BB-4627B0BEBC3F}")
